# Complain.xlsx template update
# - Adds letterhead block (UBND / Sở TN&MT, quốc hiệu, tiêu ngữ)
# - Updates report title text
# - Fixes "Địa chỉ" header typo
# - Adds signature/date line and footer app-name line

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unfreeze panes (template used to freeze header row) ---
$ws.Application.ActiveWindow.FreezePanes = $false

# --- Row 1: replace old banner text with letterhead (left) + quoc hieu (right) ---
$ws.Range("A1:D1").Merge()
$ws.Range("A1").Value = "UBND tỉnh Thái Nguyên"
$ws.Range("A1:D1").Font.Name = "Times New Roman"
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").HorizontalAlignment = -4108
$ws.Range("A1:D1").VerticalAlignment = -4108

$ws.Range("H1:J1").Merge()
$ws.Range("H1").Value = "CỘNG HÒA XÃ HỘI CHỦ NGHĨA VIỆT NAM"
$ws.Range("H1:J1").Font.Name = "Times New Roman"
$ws.Range("H1:J1").Font.Color = 0
$ws.Range("H1:J1").Font.Bold = $true
$ws.Range("H1:J1").Font.Size = 12
$ws.Range("H1:J1").HorizontalAlignment = -4108

$ws.Range("K1").Font.Name = "Times New Roman"

$ws.Rows.Item(1).RowHeight = 15.6

# --- Row 2: second letterhead line (left) + tieu ngu (right) ---
$ws.Range("A2:D2").Merge()
$ws.Range("A2").Value = "Sở Tài nguyên và Môi trường Thái Nguyên"
$ws.Range("A2:D2").Font.Name = "Times New Roman"
$ws.Range("A2:D2").Font.Bold = $true
$ws.Range("A2:D2").HorizontalAlignment = -4108
$ws.Range("A2:D2").VerticalAlignment = -4108

$ws.Range("H2:J2").Merge()
$ws.Range("H2").Value = "Độc lập - Tự do - Hạnh phúc"
$ws.Range("H2:J2").Font.Name = "Times New Roman"
$ws.Range("H2:J2").Font.Color = 0
$ws.Range("H2:J2").Font.Bold = $true
$ws.Range("H2:J2").Font.Size = 12
$ws.Range("H2:J2").HorizontalAlignment = -4108

$ws.Range("K2").Font.Name = "Times New Roman"

$ws.Rows.Item(2).RowHeight = 15.6

# --- Row 3: underline dashes under the tieu ngu ---
$ws.Range("H3:J3").Merge()
$ws.Range("H3").Value = "---------------"
$ws.Range("H3:J3").Font.Name = "Times New Roman"
$ws.Range("H3:J3").Font.Color = 0
$ws.Range("H3:J3").Font.Bold = $true
$ws.Range("H3:J3").Font.Size = 12
$ws.Range("H3:J3").HorizontalAlignment = -4108

$ws.Rows.Item(3).RowHeight = 15.6

# --- Row 4: main report title (was row 3's merged A3:F3 title) ---
$ws.Range("A3:F3").UnMerge()
$ws.Range("A3:G3").ClearContents()
$ws.Range("A3:G3").ClearFormats()
$ws.Range("A4:N4").Merge()
$ws.Range("A4").Value = "DANH SÁCH ĐƠN THƯ KHIẾU NẠI/KHIẾU KIỆN"
$ws.Range("A4:N4").Font.Name = "Times New Roman"
$ws.Range("A4:N4").Font.Bold = $true
$ws.Range("A4:N4").Font.Size = 18
$ws.Range("A4:N4").HorizontalAlignment = -4108
$ws.Rows.Item(4).RowHeight = 22.8

# --- Row 16: fix header typo "Đại chỉ" -> "Địa chỉ" + uniform header style ---
$ws.Range("D16").Value = "Địa chỉ"
$ws.Range("A16:N16").Font.Name = "Times New Roman"
$ws.Range("A16:N16").Font.Bold = $true
$ws.Range("A16:N16").WrapText = $true
$ws.Range("A16:N16").HorizontalAlignment = -4108
$ws.Range("A16:N16").VerticalAlignment = -4108
$ws.Rows.Item(16).RowHeight = 26.4

# --- Simplify the input-field boxes (remove border / date format) ---
$ws.Range("E8:E13").Borders.LineStyle = -4142
$ws.Range("E8:E13").NumberFormat = "General"

# --- New rows: signature/date line and footer app name ---
$ws.Range("I20").Value = "Thái Nguyên, Ngày …… Tháng ……. Năm ……."
$ws.Range("I20").Font.Name = "Times New Roman"
$ws.Range("I20").Font.Bold = $true

$ws.Range("B25").Value = "HTTTTT Quản lý Đơn thư Khiếu nại, Tố cáo"
$ws.Range("B25").Font.Name = "Times New Roman"

# --- View / zoom ---
$ws.Application.ActiveWindow.Zoom = 96
$ws.Range("E15").Select()
